# The deck's Slide Master design was switched from the "Integral" theme
# (Red Violet colour scheme) to the default "Office Theme" colour scheme.
# In the PowerPoint object model this is the Design > Variants > Colors
# gallery action, which rewrites each of the 12 theme colour slots
# (msoThemeColorDark1 .. msoThemeColorFollowedHyperlink) on the slide
# master's theme.

function ToOleColor([int]$r, [int]$g, [int]$b) {
    # Mirrors VBA's RGB() helper: OLE_COLOR packs bytes as 0x00BBGGRR.
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" theme colour scheme (12 slots, in
# MsoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink).
$officeColors = @(
    @(0x00,0x00,0x00),  # dk1
    @(0xFF,0xFF,0xFF),  # lt1
    @(0x44,0x54,0x6A),  # dk2
    @(0xE7,0xE6,0xE6),  # lt2
    @(0x5B,0x9B,0xD5),  # accent1
    @(0xED,0x7D,0x31),  # accent2
    @(0xA5,0xA5,0xA5),  # accent3
    @(0xFF,0xC0,0x00),  # accent4
    @(0x44,0x72,0xC4),  # accent5
    @(0x70,0xAD,0x47),  # accent6
    @(0x05,0x63,0xC1),  # hlink
    @(0x95,0x4F,0x72)   # folHlink
)

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $rgb = $officeColors[$i - 1]
    $themeColors.Item($i).RGB = ToOleColor $rgb[0] $rgb[1] $rgb[2]
}
